# Applies the 1st-iteration update to the eclaire-study-party-role-code-system
# CodeSystem spreadsheet:
#   - Metadata!B7  (Experimental)   "" -> "true"
#   - Metadata!B8  (Date)           "2025-07-17T13:20:13+00:00" -> "2025-07-21T12:46:15+00:00"
#   - Metadata!B18 (Compositional)  "" -> "false"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "true"
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "false"
